$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, shifting existing rows 164-266 down to 165-267.
$ws.Rows(164).Insert()

# Populate the newly inserted row 164 with the new record.
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 44452
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100102
$ws.Range("H164").Value = "Cítricos"
$ws.Range("I164").Value = 100102005
$ws.Range("J164").Value = "Naranja"
$ws.Range("K164").Value = "Navel Late"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 500
$ws.Range("N164").Value = 6000
$ws.Range("O164").Value = 6000
$ws.Range("P164").Value = 6000
$ws.Range("Q164").Value = "`$/caja 15 kilos granel"
$ws.Range("R164").Value = "Provincia de Quillota"
$ws.Range("S164").Value = 400
$ws.Range("T164").Value = 15
